# Scheduled market-data refresh: refresh cached Universalis price/profit
# figures (columns H:N -- currentAveragePrice*, LevePrice*, LeveProfit*)
# for the leves whose item prices moved since the last run.
#
# $wb is the already-open workbook (see harness docstring).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 1772.8182
$ws.Range("J17").Value = 1772.8182
$ws.Range("L17").Value = 5318.4546
$ws.Range("N17").Value = -5654.4546

# Row 31
$ws.Range("H31").Value = 2511.5
$ws.Range("J31").Value = 25
$ws.Range("L31").Value = 75
$ws.Range("N31").Value = -535

# Row 53
$ws.Range("H53").Value = 622.2
$ws.Range("I53").Value = 625
$ws.Range("J53").Value = 620.3333
$ws.Range("K53").Value = 625
$ws.Range("L53").Value = 620.3333
$ws.Range("M53").Value = 12
$ws.Range("N53").Value = -1894.3333

# Row 61
$ws.Range("H61").Value = 790
$ws.Range("I61").Value = 790
$ws.Range("K61").Value = 2370
$ws.Range("M61").Value = -2198

# Row 64
$ws.Range("H64").Value = 3584.4443
$ws.Range("I64").Value = 3532.75
$ws.Range("J64").Value = 3998
$ws.Range("K64").Value = 3532.75
$ws.Range("L64").Value = 3998
$ws.Range("M64").Value = -3284.75
$ws.Range("N64").Value = -4494

# Row 67
$ws.Range("H67").Value = 3584.4443
$ws.Range("I67").Value = 3532.75
$ws.Range("J67").Value = 3998
$ws.Range("K67").Value = 3532.75
$ws.Range("L67").Value = 3998
$ws.Range("M67").Value = -2674.75
$ws.Range("N67").Value = -5714

# Row 88
$ws.Range("H88").Value = 898197.4399999999
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 898197.4399999999
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 898197.4399999999
$ws.Range("M88").ClearContents()
$ws.Range("N88").Value = -899009.4399999999

# Row 91
$ws.Range("H91").Value = 898197.4399999999
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 898197.4399999999
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 898197.4399999999
$ws.Range("M91").ClearContents()
$ws.Range("N91").Value = -901005.4399999999

# Row 137
$ws.Range("H137").Value = 5201.143
$ws.Range("I137").Value = 3909
$ws.Range("J137").Value = 21999
$ws.Range("K137").Value = 11727
$ws.Range("L137").Value = 65997
$ws.Range("M137").Value = -9177
$ws.Range("N137").Value = -71097

# Row 138
$ws.Range("H138").Value = 2513.45
$ws.Range("I138").Value = 1590.1666
$ws.Range("J138").Value = 2781.5
$ws.Range("K138").Value = 4770.4998
$ws.Range("L138").Value = 8344.5
$ws.Range("M138").Value = 369.5002000000004
$ws.Range("N138").Value = -18624.5

# Row 141
$ws.Range("H141").Value = 5407
$ws.Range("I141").Value = 5597.35
$ws.Range("J141").Value = 1600
$ws.Range("K141").Value = 16792.05
$ws.Range("L141").Value = 4800
$ws.Range("M141").Value = -11612.05
$ws.Range("N141").Value = -15160

$ws = $wb.Worksheets.Item("ARM")
# Row 132
$ws.Range("H132").Value = 2759.83
$ws.Range("I132").Value = 1367.175
$ws.Range("J132").Value = 7044.923
$ws.Range("K132").Value = 4101.525
$ws.Range("L132").Value = 21134.769
$ws.Range("M132").Value = -1571.525
$ws.Range("N132").Value = -26194.769

# Row 133
$ws.Range("H133").Value = 84998.75
$ws.Range("J133").Value = 89999.336
$ws.Range("L133").Value = 89999.336
$ws.Range("N133").Value = -95059.336

$ws = $wb.Worksheets.Item("BSM")
# Row 25
$ws.Range("H25").Value = 2803.8333
$ws.Range("I25").Value = 3234.6
$ws.Range("K25").Value = 3234.6
$ws.Range("M25").Value = -2999.6

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 1121112.6
$ws.Range("I31").Value = 11762.934
$ws.Range("J31").Value = 3894486.8
$ws.Range("K31").Value = 11762.934
$ws.Range("L31").Value = 3894486.8
$ws.Range("M31").Value = -11467.934
$ws.Range("N31").Value = -3895076.8

# Row 34
$ws.Range("H34").Value = 1121112.6
$ws.Range("I34").Value = 11762.934
$ws.Range("J34").Value = 3894486.8
$ws.Range("K34").Value = 11762.934
$ws.Range("L34").Value = 3894486.8
$ws.Range("M34").Value = -11560.934
$ws.Range("N34").Value = -3894890.8

# Row 74
$ws.Range("H74").Value = 76314
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 76314
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 76314
$ws.Range("M74").ClearContents()
$ws.Range("N74").Value = -78062

# Row 77
$ws.Range("H77").Value = 76314
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 76314
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 228942
$ws.Range("M77").ClearContents()
$ws.Range("N77").Value = -237678

# Row 107
$ws.Range("H107").Value = 2229.9333
$ws.Range("I107").Value = 1032.7142
$ws.Range("K107").Value = 1032.7142
$ws.Range("M107").Value = 887.2858000000001

# Row 132
$ws.Range("H132").Value = 1986.4286
$ws.Range("I132").Value = 2036.1333
$ws.Range("J132").Value = 1688.2
$ws.Range("K132").Value = 6108.3999
$ws.Range("L132").Value = 5064.6
$ws.Range("M132").Value = -3578.3999
$ws.Range("N132").Value = -10124.6

# Row 134
$ws.Range("H134").Value = 2302.2974
$ws.Range("I134").Value = 1144.6552
$ws.Range("K134").Value = 3433.9656
$ws.Range("M134").Value = -898.9655999999995

$ws = $wb.Worksheets.Item("CUL")
# Row 31
$ws.Range("H31").Value = 940
$ws.Range("I31").Value = 980
$ws.Range("J31").Value = 900
$ws.Range("K31").Value = 2940
$ws.Range("L31").Value = 2700
$ws.Range("M31").Value = -2652
$ws.Range("N31").Value = -3276

# Row 46
$ws.Range("H46").Value = 1916
$ws.Range("I46").Value = 470.66666
$ws.Range("K46").Value = 1411.99998
$ws.Range("M46").Value = -1320.99998

# Row 61
$ws.Range("H61").Value = 2000111.8
$ws.Range("I61").Value = 2500089.8
$ws.Range("J61").Value = 200
$ws.Range("K61").Value = 7500269.399999999
$ws.Range("L61").Value = 600
$ws.Range("M61").Value = -7500054.399999999
$ws.Range("N61").Value = -1030

$ws = $wb.Worksheets.Item("GSM")
# Row 31
$ws.Range("H31").Value = 1790
$ws.Range("I31").Value = 1790
$ws.Range("K31").Value = 1790
$ws.Range("M31").Value = -1498

# Row 37
$ws.Range("H37").Value = 1790
$ws.Range("I37").Value = 1790
$ws.Range("K37").Value = 1790
$ws.Range("M37").Value = -1513

# Row 41
$ws.Range("H41").Value = 9025.5
$ws.Range("I41").Value = 10051
$ws.Range("J41").Value = 8000
$ws.Range("K41").Value = 10051
$ws.Range("L41").Value = 8000
$ws.Range("M41").Value = -9696
$ws.Range("N41").Value = -8710

# Row 70
$ws.Range("H70").Value = 6073
$ws.Range("I70").Value = 4856.2856
$ws.Range("J70").Value = 8202.25
$ws.Range("K70").Value = 4856.2856
$ws.Range("L70").Value = 8202.25
$ws.Range("M70").Value = -4586.2856
$ws.Range("N70").Value = -8742.25

# Row 73
$ws.Range("H73").Value = 6073
$ws.Range("I73").Value = 4856.2856
$ws.Range("J73").Value = 8202.25
$ws.Range("K73").Value = 4856.2856
$ws.Range("L73").Value = 8202.25
$ws.Range("M73").Value = -3920.2856
$ws.Range("N73").Value = -10074.25

# Row 74
$ws.Range("H74").Value = 40179.043
$ws.Range("J74").Value = 40179.043
$ws.Range("L74").Value = 40179.043
$ws.Range("N74").Value = -42051.043

# Row 77
$ws.Range("H77").Value = 40179.043
$ws.Range("J77").Value = 40179.043
$ws.Range("L77").Value = 120537.129
$ws.Range("N77").Value = -129897.129

# Row 107
$ws.Range("H107").Value = 453.53845
$ws.Range("I107").Value = 350.25
$ws.Range("K107").Value = 350.25
$ws.Range("M107").Value = 1569.75

# Row 113
$ws.Range("H113").Value = 5285.7144
$ws.Range("I113").Value = 6000
$ws.Range("J113").Value = 5000
$ws.Range("K113").Value = 6000
$ws.Range("L113").Value = 5000
$ws.Range("M113").Value = -3830
$ws.Range("N113").Value = -9340

# Row 122
$ws.Range("H122").Value = 2052.9092
$ws.Range("I122").Value = 2186
$ws.Range("J122").Value = 1210
$ws.Range("K122").Value = 6558
$ws.Range("L122").Value = 3630
$ws.Range("M122").Value = -4108
$ws.Range("N122").Value = -8530

# Row 132
$ws.Range("H132").Value = 14496501
$ws.Range("I132").Value = 18870088
$ws.Range("J132").Value = 8991.8125
$ws.Range("K132").Value = 56610264
$ws.Range("L132").Value = 26975.4375
$ws.Range("M132").Value = -56607734
$ws.Range("N132").Value = -32035.4375

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 44002.617
$ws.Range("I7").Value = 4736.4443
$ws.Range("J7").Value = 132351.5
$ws.Range("K7").Value = 4736.4443
$ws.Range("L7").Value = 132351.5
$ws.Range("M7").Value = -4624.4443
$ws.Range("N7").Value = -132575.5

# Row 61
$ws.Range("H61").Value = 1816.6296
$ws.Range("I61").Value = 1525
$ws.Range("K61").Value = 1525
$ws.Range("M61").Value = -1323

# Row 93
$ws.Range("H93").Value = 62503560
$ws.Range("I93").Value = 90912420
$ws.Range("J93").Value = 4080.6
$ws.Range("K93").Value = 90912420
$ws.Range("L93").Value = 4080.6
$ws.Range("M93").Value = -90911172
$ws.Range("N93").Value = -6576.6

# Row 113
$ws.Range("H113").Value = 1816.6296
$ws.Range("I113").Value = 1525
$ws.Range("K113").Value = 1525
$ws.Range("M113").Value = 645

# Row 126
$ws.Range("H126").Value = 44002.617
$ws.Range("I126").Value = 4736.4443
$ws.Range("J126").Value = 132351.5
$ws.Range("K126").Value = 14209.3329
$ws.Range("L126").Value = 397054.5
$ws.Range("M126").Value = -11739.3329
$ws.Range("N126").Value = -401994.5

$ws = $wb.Worksheets.Item("WVR")
# Row 76
$ws.Range("H76").Value = 70000
$ws.Range("J76").Value = 70000
$ws.Range("L76").Value = 70000
$ws.Range("N76").Value = -70630

# Row 79
$ws.Range("H79").Value = 70000
$ws.Range("J79").Value = 70000
$ws.Range("L79").Value = 70000
$ws.Range("N79").Value = -72184

# Row 107
$ws.Range("H107").Value = 45456760
$ws.Range("I107").Value = 71430530
$ws.Range("J107").Value = 2666.5
$ws.Range("K107").Value = 214291590
$ws.Range("L107").Value = 7999.5
$ws.Range("M107").Value = -214289670
$ws.Range("N107").Value = -11839.5

# Row 126
$ws.Range("H126").Value = 3857.6667
$ws.Range("I126").Value = 3536.75
$ws.Range("K126").Value = 10610.25
$ws.Range("M126").Value = -8140.25

# Row 132
$ws.Range("H132").Value = 2294.0195
$ws.Range("I132").Value = 1758.119
$ws.Range("J132").Value = 4794.8887
$ws.Range("K132").Value = 5274.357
$ws.Range("L132").Value = 14384.6661
$ws.Range("M132").Value = -2744.357
$ws.Range("N132").Value = -19444.6661
